$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the working-hours value for Developer2 on 09.10.2021 (row 9, column C)
$ws.Range("C9").Value = "2"

# Move the active selection to D9 (as left by the user after editing)
$null = $ws.Range("D9").Select()
